$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the line text for the second data row (D3) from "I'll be back." to "Hello Ryan"
$ws.Range("D3").Value = "Hello Ryan"

# Move the active selection to D6 (as reflected in the saved sheet view)
$ws.Range("D6").Select()
